$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 54 (quarter 01-01-2021) with revised figures ---
$ws.Range("B54").Value = 205.8
$ws.Range("C54").Value = 238.2
$ws.Range("D54").Value = 124.7
$ws.Range("E54").Value = 964.5
$ws.Range("F54").Value = 109.5
$ws.Range("G54").Value = 58.7
$ws.Range("H54").Value = 120
$ws.Range("I54").Value = 121.1
$ws.Range("J54").Value = 164.5
$ws.Range("K54").Value = 42
$ws.Range("L54").Value = 101.4
$ws.Range("M54").Value = 108.6
$ws.Range("N54").Value = 126.9
$ws.Range("O54").Value = 56.6
$ws.Range("P54").Value = 155.2
$ws.Range("Q54").Value = 90.9
$ws.Range("R54").Value = 118.3
$ws.Range("S54").Value = 147.7
$ws.Range("T54").Value = 85.6
$ws.Range("U54").Value = 102.5
$ws.Range("V54").Value = 85.5
$ws.Range("W54").Value = 94.5
$ws.Range("X54").Value = 46.4
$ws.Range("Y54").Value = 92
$ws.Range("Z54").Value = 84.7
$ws.Range("AA54").Value = 108
$ws.Range("AB54").Value = 70.6
$ws.Range("AC54").Value = 74.5
$ws.Range("AD54").Value = 74.1
$ws.Range("AE54").Value = 76.3
$ws.Range("AF54").Value = 95.7
$ws.Range("AG54").Value = 137.8
$ws.Range("AH54").Value = 85.4
$ws.Range("AI54").Value = 81.8
$ws.Range("AJ54").Value = 70.6
$ws.Range("AK54").Value = 55.2
$ws.Range("AL54").Value = 70.4
$ws.Range("AM54").Value = 115.9

# --- Add new row 55 (quarter 01-04-2021) ---
# Format as Text first so Excel stores the literal label instead of
# auto-converting the date-shaped string into a date serial number, then
# drop the temporary formatting so the cell keeps the sheet's plain style.
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "01-04-2021"
$ws.Range("A55").ClearFormats()
$ws.Range("B55").Value = 99.2
$ws.Range("C55").Value = 99
$ws.Range("D55").Value = 105.3
$ws.Range("E55").Value = 32.5
$ws.Range("F55").Value = 113.5
$ws.Range("G55").Value = 99.8
$ws.Range("H55").Value = 136.4
$ws.Range("I55").Value = 136.5
$ws.Range("J55").Value = 189
$ws.Range("K55").Value = 58.1
$ws.Range("L55").Value = 108.7
$ws.Range("M55").Value = 134.5
$ws.Range("N55").Value = 176.5
$ws.Range("O55").Value = 60.4
$ws.Range("P55").Value = 179.1
$ws.Range("Q55").Value = 101.4
$ws.Range("R55").Value = 128.2
$ws.Range("S55").Value = 131.7
$ws.Range("T55").Value = 108.2
$ws.Range("U55").Value = 133.7
$ws.Range("V55").Value = 95.6
$ws.Range("W55").Value = 108
$ws.Range("X55").Value = 41.5
$ws.Range("Y55").Value = 108.8
$ws.Range("Z55").Value = 107.5
$ws.Range("AA55").Value = 114.4
$ws.Range("AB55").Value = 81.9
$ws.Range("AC55").Value = 89.7
$ws.Range("AD55").Value = 94
$ws.Range("AE55").Value = 72.5
$ws.Range("AF55").Value = 114.6
$ws.Range("AG55").Value = 176.3
$ws.Range("AH55").Value = 99.5
$ws.Range("AI55").Value = 81.4
$ws.Range("AJ55").Value = 90.2
$ws.Range("AK55").Value = 56.6
$ws.Range("AL55").Value = 65.5
$ws.Range("AM55").Value = 121.1

Write-Output "Update applied: row54 revised, row55 added"
